# Add a date line and a contact line to the Subtitle placeholder on slide 1,
# leaving the existing "Tailoring cosmos ... deployment applications" text
# and its run structure untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shape = $s.Shapes.Item("Subtitle 2")
$tr = $shape.TextFrame.TextRange

# New paragraph: "Date: 6/28/2018"
[void]$tr.InsertAfter("`rDate: 6/28/2018")

# New paragraph: "Jordan ishii: jordanishii1@gmail.com" split across runs
[void]$tr.InsertAfter("`rJordan ")
[void]$tr.InsertAfter("ishii")
[void]$tr.InsertAfter(": jordanishii1@")
[void]$tr.InsertAfter("gmail.com")

# Trailing empty paragraph
[void]$tr.InsertAfter("`r")
